# Applies the "nuevos experimentos no convexos" update: new generator
# values for the non-convex alpha experiment (ex9.1.9 Linear, Strong
# Stationary, alpha non zero).
#
# All of the changed cells on Restricciones_del_follower / Punto_modificado /
# Vector_bf / Vector_BF hold their numbers as literal TEXT in the workbook
# (shared strings), not as numeric cells - so a plain `.Value = "..."`
# assignment would get auto-coerced to a number by Excel. We force text
# storage with NumberFormat "@" before the assignment, then reset the style
# back to Normal so the cell keeps its original (default) appearance - only
# the underlying stored type (text) needs to stick.
$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Sheet, $Cell, $Text)
    $rng = $Sheet.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# --- Restricciones_del_follower ---------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

# J_0_L0_v row (row 2)
Set-TextValue $wsFollower "A2" "-6.89670276412982 - 2x_1 + 2.3621778802767475y_1 + 0.959624318994619y_2"
Set-TextValue $wsFollower "B2" "9.39670276412982"
Set-TextValue $wsFollower "E2" "3.3000000000000003"
Set-TextValue $wsFollower "F2" "6.4"

# J_0_LP_v row (row 3)
Set-TextValue $wsFollower "A3" "2.653900531434875 + x_1 - 3x_2 - 0.468665396570741y_1 + 0.3257796049333199y_2"
Set-TextValue $wsFollower "B3" "-4.653900531434875"
Set-TextValue $wsFollower "D3" "0.7"
Set-TextValue $wsFollower "E3" "0.4"
Set-TextValue $wsFollower "F3" "7.6"

# J_Ne_L0_v row (row 4)
Set-TextValue $wsFollower "A4" "-20.705172632775163 + x_1 + x_2 + 1.6765266218790738y_1 + 2.411845315685685y_2"
Set-TextValue $wsFollower "B4" "18.065172632775162"
Set-TextValue $wsFollower "D4" "0.8"
Set-TextValue $wsFollower "E4" "10.0"

# --- Punto_modificado ---------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto "A2" "4.5"
Set-TextValue $wsPunto "B2" "1.6"
Set-TextValue $wsPunto "C2" "6.1000000000000005"
Set-TextValue $wsPunto "D2" "1.55"

# --- Vector_bf -----------------------------------------------------------
# Use the explicit tab index (5th sheet) too, for the same case-insensitive-
# name-collision reason as Vector_BF below.
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf "A2" "0.7191537150305818"
Set-TextValue $wsBf "A3" "-4.078761322236707"

# --- Vector_BF -----------------------------------------------------------
# NOTE: sheet-name lookups are case-insensitive, and this workbook has both
# "Vector_bf" and "Vector_BF" tabs - Item("Vector_BF") would resolve to the
# first (wrong) match. Use the 1-based tab index instead (6th sheet).
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF "A2" "-1.7999999999999998"
Set-TextValue $wsBF "A3" "-9.8"
Set-TextValue $wsBF "A4" "-24.872987065075712"
Set-TextValue $wsBF "A5" "-27.41552525151242"

# --- Vector_Alpha --------------------------------------------------------
# This one really is a numeric cell (no shared-string coercion needed).
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 1.71
